$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...Investidor), existem quatro..." -> "...Investidor). Existem quatro..."
#   (comma -> period, lower-case "e" -> upper-case "E")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Investidor), existem", $true, $false, $false, $false, $false,
    $true, 1, $false, "Investidor). Existem", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: drop ", sendo necessário guardar a quantidade de pesquisas que
# uma determinada palavra tem" and turn the comma after "econômico" into a
# period.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "econômico, sendo necessário guardar a quantidade de pesquisas que uma determinada palavra tem. As",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "econômico. As", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: remove the old "_GoBack" bookmark (previously sitting right
# after " e seu novo sentido (descrição)").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 4: "..., período, rendimento e a quantidade de simulações daquele
# investimento para que se possa ter uma noção dos interesses dos
# usuários." -> "..., período e rendimento."
# and re-create the "_GoBack" bookmark right before " rendimento".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ", período, rendimento e a quantidade de simulações daquele investimento para que se possa ter uma noção dos interesses dos usuários.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, ", período e rendimento.", 2) | Out-Null

$rng = $d.Content.Duplicate
$rng.Find.Execute("período e", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null
$pos = $rng.End
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
